$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: a single label cell
$ws.Range("A13").Value = "2017/4/26 alphabetCNNnetv2"

# Row 14: a single label cell
$ws.Range("A14").Value = "machine learning and casebased reasoning"

# Row 15: data row with default style
$ws.Range("A15").Value = 3
$ws.Range("B15").Value = 25
$ws.Range("C15").Value = 0.9
$ws.Range("D15").Value = "20(30)"
$ws.Range("E15").Value = 15
$ws.Range("F15").Value = 36
$ws.Range("G15").Value = 49.3

# Row 16: data row with yellow fill (same style as existing rows 8 & 11)
$ws.Range("A16").Value = 3
$ws.Range("B16").Value = 25
$ws.Range("C16").Value = 0.8
$ws.Range("D16").Value = "20(35)"
$ws.Range("E16").Value = 14
$ws.Range("F16").Value = 36
$ws.Range("G16").Value = 48.9
$ws.Range("A16:G16").Interior.Color = 65535

# Row 17: data row with new white/background fill style
$ws.Range("A17").Value = 2
$ws.Range("B17").Value = 25
$ws.Range("C17").Value = 0.8
$ws.Range("D17").Value = "30(47)"
$ws.Range("E17").Value = 21
$ws.Range("F17").Value = 36
$ws.Range("G17").Value = 113
$ws.Range("A17:G17").Interior.ThemeColor = 2
$ws.Range("A17:G17").Interior.TintAndShade = 0

# Row 18: data row with default style
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = 25
$ws.Range("C18").Value = 0.8
$ws.Range("D18").Value = "34()"
$ws.Range("E18").Value = 24
$ws.Range("F18").Value = 36
$ws.Range("G18").Value = 440.5

# Row 19 intentionally left blank

# Row 20: a single label cell
$ws.Range("A20").Value = "test"

# Row 21: data row with default style (D21 stays a plain number)
$ws.Range("A21").Value = 2
$ws.Range("B21").Value = 25
$ws.Range("C21").Value = 0.8
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 4

# Row 22: data row with yellow fill
$ws.Range("A22").Value = 3
$ws.Range("B22").Value = 25
$ws.Range("C22").Value = 0.8
$ws.Range("D22").Value = "4(12)"
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 10.8
$ws.Range("A22:G22").Interior.Color = 65535

# Update the active selection to match the edited workbook
$ws.Range("E20").Select()
